$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Set C8 to "Local" to match the weekly schedule pattern (e.g. C3, C5 already "Local")
$ws.Range("C8").Value = "Local"

# Update the active selection to C9, matching the last interacted cell
$ws.Range("C9").Select()
